$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K4").Value = 57
$ws.Range("L4").Value = 149.51
$ws.Range("K5").Value = 78
$ws.Range("L5").Value = 191.61
